# Generate Report for Handoff
#
# The "89bf7667-e4b9-450e-9eda-aff11e54eb04.md" file finished a new handoff
# round, so its handoff timestamps are refreshed across the Overview sheet
# and each locale's detail sheet (row 6 / file id 89bf7667... in each case).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column D is "Latest Handoff Date" for the 89bf7667...md row (row 6)
$wsOverview.Range("D6").Value = "2016-03-22 20:42:40"

# zh-cn detail sheet: column E is "Latest Handoff Datetime" for the 89bf7667...md row (row 6)
$wsZhCn.Range("E6").Value = "2016-03-22 20:42:37"

# de-de detail sheet: column E is "Latest Handoff Datetime" for the 89bf7667...md row (row 6)
$wsDeDe.Range("E6").Value = "2016-03-22 20:42:40"
